# carga egreso, infractorcontraventor y titular
# Replace the numeric idtipodni codes in column K with their textual labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$map = @{
    2  = "DNI"
    3  = "LE"
    4  = "LC"
    5  = "LE"
    6  = "DNI"
    7  = "EXT"
    8  = "LE"
    9  = "LC"
    10 = "DNI"
    11 = "EXT"
    12 = "LE"
    13 = "LC"
    14 = "DNI"
    15 = "EXT"
    16 = "LE"
    17 = "LE"
    18 = "LC"
    19 = "DNI"
    20 = "DNI"
    21 = "DNI"
}

foreach ($row in $map.Keys) {
    $ws.Cells.Item($row, 11).Value = $map[$row]
}

$ws.Range("J24").Select()
